$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current.AddDays(1)
    }
}
